$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (pushes existing data rows down by one)
$ws.Rows.Item(2).Insert()

# New first data row: latest date (43949) and death count (195)
$ws.Range("A2").Value = 43949
$ws.Range("A2").Style = $ws.Range("A3").Style
$ws.Cells.Item(2, 1).NumberFormat = "yyyy-mm-dd"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = " 195"
$ws.Range("B2").Style = "Normal"
